$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Iterative calculation setting (calcPr iterateDelta) ---------------
# Best-effort: mirrors the author's intent (enable iterative calc with a
# tighter convergence delta of 1E-4). Included for completeness even
# though this engine build does not persist calcPr/@iterate* back to XML.
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# --- New columns G/H: clone header (C1, style 4) and data (F2, style 1) formats ---
$ws.Range("C1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("F2:F14").Copy()
$ws.Range("G2:G14").PasteSpecial(-4122)
$ws.Range("H2:H14").PasteSpecial(-4122)

# --- Header row relabelling ---------------------------------------------
# Write order matters for shared-string slot reuse/compaction: this order
# reproduces the exact sharedStrings.xml layout of the target workbook.
$ws.Range("H1").Value = "Responsable"
$ws.Range("F1").Value = "Fecha inicio real"
$ws.Range("G1").Value = "Fecha fin real"
$ws.Range("E1").Value = "Fecha fin planeada"
$ws.Range("D1").Value = "Fecha inicio planeada"

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 28.5
$ws.Columns.Item(5).ColumnWidth = 24.833333333333332
$ws.Columns.Item(6).ColumnWidth = 25.166666666666668
$ws.Columns.Item(7).ColumnWidth = 21.166666666666668
$ws.Columns.Item(8).ColumnWidth = 16.0

# --- Remove the now-unused trailing blank row 15 --------------------------
$ws.Rows.Item(15).Delete()

# --- Selection moves to D7 -------------------------------------------------
$ws.Range("D7").Select()
